$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.692.67'
$ws.Range('E2').Value = '  -2.45%  '
$ws.Range('D3').Value = '1.875.97'
$ws.Range('E3').Value = '  -2.03%  '
$ws.Range('E4').Value = '  -0.81%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '247.30'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('E6').Value = '  -3.15%  '
$ws.Range('E7').Value = '  -0.89%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.19'
$ws.Range('E8').Value = '  +1.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.346'
$ws.Range('E9').Value = '  -2.69%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '50.73'
$ws.Range('E10').Value = '  -4.65%  '
$ws.Range('E11').Value = '  +0.65%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0966'
$ws.Range('E12').Value = '  -2.41%  '
$ws.Range('D13').Value = '2.147.74'
$ws.Range('E13').Value = '  -2.04%  '
$ws.Range('E14').Value = '  +1.23%  '
$ws.Range('E15').Value = '  -0.64%  '
$ws.Range('E16').Value = '  -0.73%  '
$ws.Range('D17').Value = '1.868.77'
$ws.Range('E17').Value = '  -2.40%  '
$ws.Range('D18').Value = '34.699.78'
$ws.Range('E18').Value = '  -2.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.81'
$ws.Range('E19').Value = '  -0.81%  '
$ws.Range('D20').Value = '0.0₃0822'
$ws.Range('E20').Value = '  -0.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '246.57'
$ws.Range('E21').Value = '  +1.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.70'
$ws.Range('E22').Value = '  -3.65%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.91'
$ws.Range('E23').Value = '  -3.18%  '
$ws.Range('E24').Value = '  -0.92%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.40'
$ws.Range('E25').Value = '  +3.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.24'
$ws.Range('E26').Value = '  -2.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.38'
$ws.Range('E27').Value = '  -1.73%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.37'
$ws.Range('E28').Value = '  -3.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.23'
$ws.Range('E29').Value = '  -3.39%  '
$ws.Range('E30').Value = '  -4.04%  '
$ws.Range('D31').Value = '4.128.39'
$ws.Range('E32').Value = '  +13.61%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.23'
$ws.Range('E33').Value = '  -0.99%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0579'
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.14'
$ws.Range('E36').Value = '  -0.90%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.82'
$ws.Range('E37').Value = '  -5.97%  '
$ws.Range('E38').Value = '  -9.36%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.97'
$ws.Range('E39').Value = '  -4.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.45'
$ws.Range('E40').Value = '  -3.19%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.12'
$ws.Range('E41').Value = '  -0.92%  '
$ws.Range('E42').Value = '  +1.87%  '
$ws.Range('E43').Value = '  -0.64%  '
$ws.Range('E44').Value = '  -6.00%  '
$ws.Range('D45').Value = '1.291.87'
$ws.Range('E45').Value = '  -4.41%  '
$ws.Range('E46').Value = '  -4.98%  '
$ws.Range('E47').Value = '  -0.91%  '
$ws.Range('E48').Value = '  -1.88%  '
$ws.Range('B49').Value = 'Gas'
$ws.Range('C49').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '12.22'
$ws.Range('E49').Value = '  +0.56%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0765'
$ws.Range('E50').Value = '  +5.92%  '
$ws.Range('E51').Value = '  -1.37%  '
